# Populate additional product rows (4-23) with new names/prices pulled in
# from the refreshed catalog export. Rows 1-3 (existing products) and the
# trailing blank rows 24-30 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 4;  A = "Claire Rectangular Solid Wood Coffee Table In Teak Finish";      B = "₹12,725" },
    @{ Row = 5;  A = "Tate Square Solid Wood Coffee Table In Teak Finish";             B = "₹16,554" },
    @{ Row = 6;  A = "Striado Rectangular Solid Wood Coffee Table In Teak Finish";     B = "₹10,947" },
    @{ Row = 7;  A = "Marcel Rectangular Metal Coffee Table In White Gloss Finish";    B = "₹11,967" },
    @{ Row = 8;  A = "Renesme Rectangular Solid Wood Coffee Table In Mahogany Finish"; B = "₹15,317" },
    @{ Row = 9;  A = "Dyson Abstract Metal Coffee Table In Teak Finish";               B = "₹7,679" },
    @{ Row = 10; A = "Ivara Rectangular Solid Wood Coffee Table In Natural Finish";    B = "₹16,049" },
    @{ Row = 11; A = "Botwin Rectangular Solid Wood Coffee Table In Mahogany Finish";  B = "₹9,647" },
    @{ Row = 12; A = "Zephyr Rectangular Solid Wood Coffee Table In Teak Finish";      B = "₹14,104" },
    @{ Row = 13; A = "Fring Engineered Wood Side Table In Matte Finish";               B = "₹2,399" },
    @{ Row = 14; A = "Claire Rectangular Solid Wood Coffee Table In Mahogany Finish";  B = "₹12,725" },
    @{ Row = 15; A = "Botwin Rectangular Solid Wood Coffee Table In Teak Finish";      B = "₹9,647" },
    @{ Row = 16; A = "Epsilon Rectangular Solid Wood Coffee Table In Mahogany Finish"; B = "₹11,384" },
    @{ Row = 17; A = "Dyson Rectangular Metal Coffee Table In Walnut Finish";          B = "₹10,529" },
    @{ Row = 18; A = "Gustowe Rectangular Engineered Wood Coffee Table In Matte Finish"; B = "₹2,279" },
    @{ Row = 19; A = "Striado Rectangular Solid Wood Coffee Table In Mahogany Finish"; B = "₹10,947" },
    @{ Row = 20; A = "Osiris Rectangular Stone Coffee Table In Finish";                B = "₹15,677" },
    @{ Row = 21; A = "Altura Rectangular Solid Wood Coffee Table In Two Tone Finish";  B = "₹8,374" },
    @{ Row = 22; A = "Sylvie Rectangular Solid Wood Coffee Table In Natural Finish";   B = "₹11,839" },
    @{ Row = 23; A = "Florence Oval Solid Wood Coffee Table In Teak Finish";           B = "₹10,223" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}
